# Estado de Cuenta (NIT 9004223571) update
# - Refresh the "VALOR MORA" (overdue amount) total and the "Cant. Periodos"
#   (period count) now that only one overdue period remains.
# - Remove the data row for period 2503 (it is no longer outstanding),
#   leaving only period 2504. All rows below shift up automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated totals: Valor Mora and Cant. Periodos now reflect a single
# outstanding period (2504) instead of two (2504 + 2503).
$ws.Range("E11").Value2 = 57200
$ws.Range("F13").Value2 = 1

# Drop the whole worksheet row that held the now-settled period 2503
# (B17:J17). Excel shifts the remaining rows (signature block, etc.) up.
$ws.Rows("17").Delete()
